# Adding MPA test automation upload file
# Update transfer/upload identifier values on the "Data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Column K: 60000215 -> 60000231
foreach ($r in @(6, 9, 11, 14, 16, 19, 21, 24, 26, 29)) {
    $ws.Cells.Item($r, 11).Value = 60000231
}

# Column L: 165 -> 179
foreach ($r in @(7, 8, 10, 12, 13, 15, 17, 18, 20, 22, 23, 25, 27, 28)) {
    $ws.Cells.Item($r, 12).Value = 179
}

# Column N: 60000216 -> 60000232
foreach ($r in @(7, 11, 12, 16, 17, 21, 22, 26, 27)) {
    $ws.Cells.Item($r, 14).Value = 60000232
}

# Column O: 166 -> 180
foreach ($r in @(8, 13, 18, 23, 28)) {
    $ws.Cells.Item($r, 15).Value = 180
}
